$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" ----
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-21 16:22:21"

$ws1.Range("A3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-21 16:22:21"

$ws1.Range("A4").Value = "ed807c42-2548-4889-a845-e2dc2186f24d.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"
$ws1.Range("D4").Value = "2016-03-21 16:20:26"

$ws1.Range("A5").Value = "c89b8ba0-486a-43cd-b669-9289af8e184f.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-21 16:22:21"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/5234f991-bf43-4261-b693-6cb48c697bd4.md", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/fea4c374-084f-4455-9e46-f8d56cacb6fc.md", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bc695e88f97c5bac7f0c7eee1a8d22374d38c4db/e2e/ed807c42-2548-4889-a845-e2dc2186f24d.md", "", "", "ed807c42-2548-4889-a845-e2dc2186f24d.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62a6cf7326a9e59f5f22e04f0563c91f99217fc/e2e/c89b8ba0-486a-43cd-b669-9289af8e184f.md", "", "", "c89b8ba0-486a-43cd-b669-9289af8e184f.md")

# ---- Sheet "zh-cn" ----
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-21 16:22:17"
$ws2.Range("F2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.md"
$ws2.Range("G2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-21 16:22:38"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-21 16:22:17"
$ws2.Range("F3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.md"
$ws2.Range("G3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-21 16:22:38"
$ws2.Range("J3").Value = "Include"

$ws2.Range("A4").Value = "ed807c42-2548-4889-a845-e2dc2186f24d.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "In Translation"
$ws2.Range("D4").Value = "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-21 16:20:23"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"

$ws2.Range("A5").Value = "c89b8ba0-486a-43cd-b669-9289af8e184f.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-21 16:22:17"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/5234f991-bf43-4261-b693-6cb48c697bd4.md", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67be46802f7734564b3931c91a34941f15d719cf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/5234f991-bf43-4261-b693-6cb48c697bd4.md", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67be46802f7734564b3931c91a34941f15d719cf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/fea4c374-084f-4455-9e46-f8d56cacb6fc.md", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67be46802f7734564b3931c91a34941f15d719cf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/fea4c374-084f-4455-9e46-f8d56cacb6fc.md", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67be46802f7734564b3931c91a34941f15d719cf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bc695e88f97c5bac7f0c7eee1a8d22374d38c4db/e2e/ed807c42-2548-4889-a845-e2dc2186f24d.md", "", "", "ed807c42-2548-4889-a845-e2dc2186f24d.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10539a357a8801ade2f5e4812b9ead109d5d77d3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.zh-cn.xlf", "", "", "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62a6cf7326a9e59f5f22e04f0563c91f99217fc/e2e/c89b8ba0-486a-43cd-b669-9289af8e184f.md", "", "", "c89b8ba0-486a-43cd-b669-9289af8e184f.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67be46802f7734564b3931c91a34941f15d719cf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.zh-cn.xlf", "", "", "c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.zh-cn.xlf")

# ---- Sheet "de-de" ----
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-21 16:22:21"
$ws3.Range("F2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.md"
$ws3.Range("G2").Value = "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-21 16:22:44"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-21 16:22:21"
$ws3.Range("F3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.md"
$ws3.Range("G3").Value = "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-21 16:22:44"
$ws3.Range("J3").Value = "Include"

$ws3.Range("A4").Value = "ed807c42-2548-4889-a845-e2dc2186f24d.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "In Translation"
$ws3.Range("D4").Value = "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-21 16:20:26"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"

$ws3.Range("A5").Value = "c89b8ba0-486a-43cd-b669-9289af8e184f.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-21 16:22:21"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/5234f991-bf43-4261-b693-6cb48c697bd4.md", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a7689768e4bc98a12c4155523fcd388a0ce2787/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/5234f991-bf43-4261-b693-6cb48c697bd4.md", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a7689768e4bc98a12c4155523fcd388a0ce2787/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf", "", "", "5234f991-bf43-4261-b693-6cb48c697bd4.2f4a76a3f231122e90e8f695066ae61881c56181.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/fea4c374-084f-4455-9e46-f8d56cacb6fc.md", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a7689768e4bc98a12c4155523fcd388a0ce2787/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/1381196c372cdaccbe14aec5f25ed9388f6d6a1e/e2e/fea4c374-084f-4455-9e46-f8d56cacb6fc.md", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a7689768e4bc98a12c4155523fcd388a0ce2787/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf", "", "", "fea4c374-084f-4455-9e46-f8d56cacb6fc.81b042d97e79ac03eddd2bfb1c2c5d5155e584a9.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bc695e88f97c5bac7f0c7eee1a8d22374d38c4db/e2e/ed807c42-2548-4889-a845-e2dc2186f24d.md", "", "", "ed807c42-2548-4889-a845-e2dc2186f24d.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a5813561d9ddc0a8050a584d419313ff839a5f7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.de-de.xlf", "", "", "ed807c42-2548-4889-a845-e2dc2186f24d.a729299b0fe48fc16a80d1be6bd2e45fcf7bc22c.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d62a6cf7326a9e59f5f22e04f0563c91f99217fc/e2e/c89b8ba0-486a-43cd-b669-9289af8e184f.md", "", "", "c89b8ba0-486a-43cd-b669-9289af8e184f.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a7689768e4bc98a12c4155523fcd388a0ce2787/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.de-de.xlf", "", "", "c89b8ba0-486a-43cd-b669-9289af8e184f.c9f8da61a8b7fda2b5f3526c0a2cbba5195a0da5.de-de.xlf")
